# The document contains six inline pictures (drawings). Word normally marks
# the run that hosts a picture as "no proofing" (<w:noProof/> in w:rPr) so
# that spelling/grammar check skips image runs. This document was missing
# that flag on every picture run ("fehlende basic logiken erggänzt" - add
# the missing basic/no-proofing logic back in). Re-apply it to every
# inline picture's run.
$d = $word.ActiveDocument

foreach ($ishp in $d.InlineShapes) {
    $ishp.Range.NoProofing = 1
}
